$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the 31cc10a6 file
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-17 02:51:06"

# zh-cn sheet: row 3 is the 31cc10a6 file
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-17 02:50:58"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83b8e7ba0d7eea642ec2fb1f50bafb520ee736e7/e2e/31cc10a6-7f3a-4cc1-9095-d4ff8f021ccf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d361a527c19c24961e00dcbb2ca90f34df6c7fe1/e2e/31cc10a6-7f3a-4cc1-9095-d4ff8f021ccf.md."

# de-de sheet: row 3 is the 31cc10a6 file
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-17 02:51:06"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83b8e7ba0d7eea642ec2fb1f50bafb520ee736e7/e2e/31cc10a6-7f3a-4cc1-9095-d4ff8f021ccf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d361a527c19c24961e00dcbb2ca90f34df6c7fe1/e2e/31cc10a6-7f3a-4cc1-9095-d4ff8f021ccf.md."

# Widen column P (Error Detail) on zh-cn and de-de sheets
# (39.17 round-trips to exactly width="40" in the saved OOXML)
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
